{"js": "const body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\n// Locate the five \"ORDER:\" list paragraphs by their original text so the\n// script is resilient to exact indices.\nconst targets = {\n  playerInventory: null,\n  gainAether: null,\n  stageProgression: null,\n  stageSelectUi: null,\n};\nfor (let i = 0; i < paras.items.length; i++) {\n  const t = paras.items[i].text;\n  if (t === \"Player inventory\" && targets.playerInventory === null) {\n    targets.playerInventory = i;\n  } else if (t === \"Gain aether from finishing a stage\" && targets.gainAether === null) {\n    targets.gainAether = i;\n  } else if (t === \"Stage progression\" && targets.stageProgression === null) {\n    targets.stageProgression = i;\n  } else if (t === \"Stage select ui\" && targets.stageSelectUi === null) {\n    targets.stageSelectUi = i;\n  }\n}\n\n// 1) \"Player inventory\" -> two runs of brainstorming text in the same paragraph.\nconst pPlayerInventory = paras.items[targets.playerInventory];\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t>Look back over uiElement and uiState.  Is it actually the way I want with the overrides/post constructor?  Is there a better way?</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">  I could maybe make them separate objects and have common methods between them instead of inheriting from each other.  AKA, all inherit from uiElement, but have different implementations for rect?</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\npPlayerInventory.getRange(\"Whole\").insertOoxml(ooxml, \"Replace\");\n\n// 2) \"Gain aether from finishing a stage\" -> \"Side bar UI\"\nparas.items[targets.gainAether].insertText(\"Side bar UI\", \"Replace\");\n\n// 3) \"Stage progression\" -> \"Animations and UIState\"\nparas.items[targets.stageProgression].insertText(\"Animations and UIState\", \"Replace\");\n\n// 4) New paragraph \"Stage select ui\" inserted right after (former) \"Stage progression\".\nparas.items[targets.stageProgression].insertParagraph(\"Stage select ui\", \"After\");\n\nawait context.sync();\n\n// 5) The trailing (second) \"Stage select ui\" paragraph becomes \"Stage data saving\".\n// \"Stage bonus manager\" stays untouched in between, so re-load paragraphs and\n// find the LAST paragraph whose text is \"Stage select ui\" (the one that was\n// already present before our edits, now pushed down by the insertion above).\nconst paras2 = body.paragraphs;\nparas2.load(\"items/text\");\nawait context.sync();\n\nlet lastStageSelectUi = -1;\nfor (let i = 0; i < paras2.items.length; i++) {\n  if (paras2.items[i].text === \"Stage select ui\") {\n    lastStageSelectUi = i;\n  }\n}\nparas2.items[lastStageSelectUi].insertText(\"Stage data saving\", \"Replace\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Find-ExactParagraphIndex($doc, $text) {\n  $n = $doc.Paragraphs.Count\n  for ($i = 1; $i -le $n; $i++) {\n    $p = $doc.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    $t2 = $t.TrimEnd([char]13)\n    if ($t2 -eq $text) {\n      return $i\n    }\n  }\n  return -1\n}\n\n# 1) \"Stage progression\" -> \"Animations and UIState\"\n$idx = Find-ExactParagraphIndex $d \"Stage progression\"\n$p = $d.Paragraphs.Item($idx)\n$r = $p.Range\n[void]$r.MoveEnd(1, -1)\n$r.Text = \"Animations and UIState\"\n\n# 2) Insert a brand-new paragraph \"Stage select ui\" right after the paragraph\n#    we just renamed (was \"Stage progression\").\n$p = $d.Paragraphs.Item($idx)\n$r = $p.Range\n[void]$r.InsertParagraphAfter()\n$newIdx = $idx + 1\n$newP = $d.Paragraphs.Item($newIdx)\n$nr = $newP.Range\n[void]$nr.MoveEnd(1, -1)\n$nr.Text = \"Stage select ui\"\n\n# 3) \"Gain aether from finishing a stage\" -> \"Side bar UI\"\n$idx = Find-ExactParagraphIndex $d \"Gain aether from finishing a stage\"\n$p = $d.Paragraphs.Item($idx)\n$r = $p.Range\n[void]$r.MoveEnd(1, -1)\n$r.Text = \"Side bar UI\"\n\n# 4) \"Player inventory\" -> a single paragraph containing TWO separate runs of\n#    brainstorming text.\n$idx = Find-ExactParagraphIndex $d \"Player inventory\"\n$p = $d.Paragraphs.Item($idx)\n$r = $p.Range\n[void]$r.MoveEnd(1, -1)\n$r.Text = \"\"\n$ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t>Look back over uiElement and uiState.  Is it actually the way I want with the overrides/post constructor?  Is there a better way?</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">  I could maybe make them separate objects and have common methods between them instead of inheriting from each other.  AKA, all inherit from uiElement, but have different implementations for rect?</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>'\n[void]$r.InsertXML($ooxml)\n\n# 5) The (still remaining, original) \"Stage select ui\" paragraph that used to\n#    sit at the end of the list becomes \"Stage data saving\". \"Stage bonus\n#    manager\" sits between the two \"Stage select ui\" paragraphs and must stay\n#    untouched, so take the LAST paragraph matching \"Stage select ui\".\n$n = $d.Paragraphs.Count\n$lastIdx = -1\nfor ($i = 1; $i -le $n; $i++) {\n  $t = $d.Paragraphs.Item($i).Range.Text\n  $t2 = $t.TrimEnd([char]13)\n  if ($t2 -eq \"Stage select ui\") {\n    $lastIdx = $i\n  }\n}\n$p = $d.Paragraphs.Item($lastIdx)\n$r = $p.Range\n[void]$r.MoveEnd(1, -1)\n$r.Text = \"Stage data saving\"\n\nWrite-Output \"done\"\n"}
